$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column K (new "miles" helper column) to roughly match the authored layout.
$ws.Columns("K").ColumnWidth = 17.9

# Row 27 - new mini "km -> miles" table headers
$ws.Range("J27").Value = "km"
$ws.Range("K27").Value = "fun"
$ws.Range("L27").Value = "miles"

# Row 28 - first data row; L28 is a plain literal value (not a formula)
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = "miles=0.621371*km"
$ws.Range("L28").Value = 0.621371

# Rows 29-34 - remaining data rows with J = km input, L = formula conversion
$ws.Range("J29").Value = 2
$ws.Range("J30").Value = 3
$ws.Range("J31").Value = 4
$ws.Range("J32").Value = 5
$ws.Range("J33").Value = 6
$ws.Range("J34").Value = 25

$ws.Range("L29").Formula = "=J29*0.621371"
$ws.Range("L30").Formula = "=J30*0.621371"
$ws.Range("L31").Formula = "=J31*0.621371"
$ws.Range("L32").Formula = "=J32*0.621371"
$ws.Range("L33").Formula = "=J33*0.621371"
$ws.Range("L34").Formula = "=J34*0.621371"

# Move the selection/view to where the new table was added.
$ws.Range("H33").Select() | Out-Null
